# Restructuration of transitions matrix.
#
# The worksheet "Hoja1" holds a lexer transition matrix: row = current
# state (col A), columns B..AO = input characters/char-classes, and each
# intersection cell holds either the literal "ER" (no transition / error)
# or the numeric id of the destination state. This script moves several
# transition values to new (character, state) pairs, fills in previously
# blank rows 27/28 with the full "ER" grid plus new destination states,
# and tags two procedure markers (PR01/PR02) used by rows 18/28.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9 (state 0) ---------------------------------------------------
$ws.Range("F9").Value = 20
$ws.Range("G9").Value = 28
$ws.Range("Q9").Value = 10
$ws.Range("T9").Value = 41
$ws.Range("AN9").Value = 1

# --- Row 10 (state 1) ---------------------------------------------------
$ws.Range("F10").Value = 20
$ws.Range("G10").Value = 28
$ws.Range("Q10").Value = 10
$ws.Range("T10").Value = 41
$ws.Range("AL10").Value = 1
$ws.Range("AN10").Value = 1

# --- Row 17 (state 8) ---------------------------------------------------
$ws.Range("AL17").Value = 9
$ws.Range("AN17").Value = 9

# --- Row 18 (state 9) ---------------------------------------------------
$ws.Range("S18").Value = "ER"
# AM18 keeps its "ER" value but is now rendered underlined.
$ws.Range("AM18").Font.Underline = $true
$ws.Range("AP18").Value = "PR01"
$ws.Range("AQ18").Value = "incluir"

# --- Row 19 (state 10) ---------------------------------------------------
$ws.Range("J19").Value = "ER"
$ws.Range("S19").Value = 11

# --- Row 20 (state 11) ---------------------------------------------------
$ws.Range("J20").Value = 12
$ws.Range("O20").Value = "ER"

# --- Row 21 (state 12) ---------------------------------------------------
$ws.Range("D21").Value = "ER"
$ws.Range("O21").Value = 13

# --- Row 22 (state 13) ---------------------------------------------------
$ws.Range("D22").Value = 14
$ws.Range("J22").Value = "ER"

# --- Row 23 (state 14) ---------------------------------------------------
$ws.Range("J23").Value = 15
$ws.Range("Q23").Value = "ER"

# --- Row 24 (state 15) ---------------------------------------------------
$ws.Range("B24").Value = "ER"
$ws.Range("Q24").Value = 16

# --- Row 25 (state 16) ---------------------------------------------------
$ws.Range("B25").Value = 17
$ws.Range("M25").Value = "ER"

# --- Row 26 (state 17) ---------------------------------------------------
$ws.Range("M26").Value = 18

# --- Row 27 (state 18) -- previously blank, now a full "ER" row with a
#     new destination state 19 in the whitespace columns (AL/AM/AN). ----
for ($col = 2; $col -le 37; $col++) {
    $ws.Cells.Item(27, $col).Value = "ER"
}
$ws.Cells.Item(27, 38).Value = 19   # AL27 (" ")
$ws.Cells.Item(27, 39).Value = 19   # AM27 ("\r")
$ws.Cells.Item(27, 40).Value = 19   # AN27 ("\t")
$ws.Cells.Item(27, 41).Value = "ER" # AO27 ("_")

# --- Row 28 (state 19) -- previously blank, now a full "ER" row plus the
#     PR02/principal markers in the two new trailing columns. -----------
for ($col = 2; $col -le 41; $col++) {
    $ws.Cells.Item(28, $col).Value = "ER"
}
$ws.Cells.Item(28, 42).Value = "PR02"     # AP28
$ws.Cells.Item(28, 43).Value = "principal" # AQ28

# Restore the cursor/selection to AM18, matching the saved view state.
$ws.Range("AM18").Select() | Out-Null
